$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.193982601165771
$ws.Range("B1").Value = 6.009495258331299
$ws.Range("C1").Value = 4.91078519821167
$ws.Range("D1").Value = 5.704831123352051
$ws.Range("E1").Value = 5.061182975769043
